$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28 (anchor G28=27772)
$ws.Range("H28").Value = 848.1667
$ws.Range("J28").Value = 963
$ws.Range("L28").Value = 963
$ws.Range("N28").Value = -1933
# Row 137 (anchor G137=44013)
$ws.Range("H137").Value = 4287.4443
$ws.Range("I137").Value = 5347.8335
$ws.Range("J137").Value = 2166.6667
$ws.Range("K137").Value = 16043.5005
$ws.Range("L137").Value = 6500.000100000001
$ws.Range("M137").Value = -13493.5005
$ws.Range("N137").Value = -11600.0001
# Row 138 (anchor G138=44169)
$ws.Range("H138").Value = 2728.652
$ws.Range("I138").Value = 3178
$ws.Range("J138").Value = 2532.0625
$ws.Range("K138").Value = 9534
$ws.Range("L138").Value = 7596.1875
$ws.Range("M138").Value = -4394
$ws.Range("N138").Value = -17876.1875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45 (anchor G45=27714)
$ws.Range("H45").Value = 7593.5557
$ws.Range("I45").Value = 2974.75
$ws.Range("K45").Value = 2974.75
$ws.Range("M45").Value = -2597.75
# Row 74 (anchor G74=44000)
$ws.Range("H74").Value = 61342.785
$ws.Range("I74").Value = 103005.375
$ws.Range("J74").Value = 5792.6665
$ws.Range("K74").Value = 103005.375
$ws.Range("L74").Value = 5792.6665
$ws.Range("M74").Value = -102131.375
$ws.Range("N74").Value = -7540.6665
# Row 77 (anchor G77=44000)
$ws.Range("H77").Value = 61342.785
$ws.Range("I77").Value = 103005.375
$ws.Range("J77").Value = 5792.6665
$ws.Range("K77").Value = 515026.875
$ws.Range("L77").Value = 28963.3325
$ws.Range("M77").Value = -510658.875
$ws.Range("N77").Value = -37699.3325

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86 (anchor G86=12526)
$ws.Range("H86").Value = 35754044
$ws.Range("I86").Value = 50915.523
$ws.Range("J86").Value = 142863420
$ws.Range("K86").Value = 50915.523
$ws.Range("L86").Value = 142863420
$ws.Range("M86").Value = -49792.523
$ws.Range("N86").Value = -142865666
# Row 89 (anchor G89=12526)
$ws.Range("H89").Value = 35754044
$ws.Range("I89").Value = 50915.523
$ws.Range("J89").Value = 142863420
$ws.Range("K89").Value = 254577.615
$ws.Range("L89").Value = 714317100
$ws.Range("M89").Value = -248961.615
$ws.Range("N89").Value = -714328332
# Row 104 (anchor G104=19571)
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = $null
# Row 105 (anchor G105=19947)
$ws.Range("H105").Value = 79239.10000000001
$ws.Range("I105").Value = 92322.82000000001
$ws.Range("J105").Value = 5098
$ws.Range("K105").Value = 92322.82000000001
$ws.Range("L105").Value = 5098
$ws.Range("M105").Value = -90575.82000000001
$ws.Range("N105").Value = -8592
# Row 134 (anchor G134=43998)
$ws.Range("H134").Value = 6574.757
$ws.Range("J134").Value = 8911.708000000001
$ws.Range("L134").Value = 26735.124
$ws.Range("N134").Value = -31805.124

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (anchor G16=27691)
$ws.Range("H16").Value = 2875.516
$ws.Range("I16").Value = 1437.1578
$ws.Range("J16").Value = 5152.9165
$ws.Range("K16").Value = 1437.1578
$ws.Range("L16").Value = 5152.9165
$ws.Range("M16").Value = -1150.1578
$ws.Range("N16").Value = -5726.9165
# Row 112 (anchor G112=25796)
$ws.Range("H112").Value = 78000
$ws.Range("J112").Value = 78000
$ws.Range("L112").Value = 78000
$ws.Range("N112").Value = -80954
# Row 113 (anchor G113=27691)
$ws.Range("H113").Value = 2875.516
$ws.Range("I113").Value = 1437.1578
$ws.Range("J113").Value = 5152.9165
$ws.Range("K113").Value = 1437.1578
$ws.Range("L113").Value = 5152.9165
$ws.Range("M113").Value = 732.8422
$ws.Range("N113").Value = -9492.916499999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 99 (anchor G99=19817)
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = $null
# Row 100 (anchor G100=19831)
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = $null
# Row 101 (anchor G101=19820)
$ws.Range("H101").Value = 2000
$ws.Range("J101").Value = 2000
$ws.Range("L101").Value = 6000
$ws.Range("N101").Value = -10868
# Row 102 (anchor G102=19813)
$ws.Range("H102").Value = 15000
$ws.Range("J102").Value = 15000
$ws.Range("L102").Value = 45000
$ws.Range("N102").Value = -49868
# Row 103 (anchor G103=19839)
$ws.Range("H103").Value = 1703.8572
$ws.Range("I103").Value = 900
$ws.Range("J103").Value = 1837.8334
$ws.Range("K103").Value = 2700
$ws.Range("L103").Value = 5513.5002
$ws.Range("M103").Value = -1821
$ws.Range("N103").Value = -7271.5002
# Row 104 (anchor G104=19807)
$ws.Range("H104").Value = 3141
$ws.Range("I104").Value = 4665.3335
$ws.Range("J104").Value = 1997.75
$ws.Range("K104").Value = 13996.0005
$ws.Range("L104").Value = 5993.25
$ws.Range("M104").Value = -11375.0005
$ws.Range("N104").Value = -11235.25
# Row 105 (anchor G105=19814)
$ws.Range("H105").Value = 12000
$ws.Range("J105").Value = 12000
$ws.Range("L105").Value = 36000
$ws.Range("N105").Value = -41242

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 32 (anchor G32=27215)
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = $null
# Row 54 (anchor G54=2130)
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = $null
# Row 97 (anchor G97=19940)
$ws.Range("H97").Value = 907.13513
$ws.Range("I97").Value = 847.6875
$ws.Range("K97").Value = 847.6875
$ws.Range("M97").Value = -351.6875
# Row 99 (anchor G99=19532)
$ws.Range("H99").Value = 5143.75
$ws.Range("I99").Value = 1400
$ws.Range("J99").Value = 8887.5
$ws.Range("K99").Value = 1400
$ws.Range("L99").Value = 8887.5
$ws.Range("M99").Value = 846
$ws.Range("N99").Value = -13379.5
# Row 118 (anchor G118=26172)
$ws.Range("H118").Value = 15999.5
$ws.Range("J118").Value = 15999.5
$ws.Range("L118").Value = 15999.5
$ws.Range("N118").Value = -19313.5
# Row 121 (anchor G121=26338)
$ws.Range("H121").Value = 50030.668
$ws.Range("J121").Value = 50030.668
$ws.Range("L121").Value = 50030.668
$ws.Range("N121").Value = -53524.668
# Row 132 (anchor G132=44008)
$ws.Range("H132").Value = 4237.64
$ws.Range("I132").Value = 1694.375
$ws.Range("K132").Value = 5083.125
$ws.Range("M132").Value = -2553.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (anchor G22=5277)
$ws.Range("H22").Value = 1893.4
$ws.Range("I22").Value = 955.55554
$ws.Range("K22").Value = 955.55554
$ws.Range("M22").Value = -660.55554
# Row 27 (anchor G27=5277)
$ws.Range("H27").Value = 1893.4
$ws.Range("I27").Value = 955.55554
$ws.Range("K27").Value = 955.55554
$ws.Range("M27").Value = -848.55554
# Row 61 (anchor G61=27740)
$ws.Range("H61").Value = 6842.778
$ws.Range("I61").Value = 4663.8887
$ws.Range("K61").Value = 4663.8887
$ws.Range("M61").Value = -4461.8887
# Row 68 (anchor G68=12563)
$ws.Range("H68").Value = 3611.5454
$ws.Range("I68").Value = 2191.889
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 2191.889
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -1442.889
$ws.Range("N68").Value = -11498
# Row 71 (anchor G71=12563)
$ws.Range("H71").Value = 3611.5454
$ws.Range("I71").Value = 2191.889
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 10959.445
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -7215.445
$ws.Range("N71").Value = -57488
# Row 95 (anchor G95=18221)
$ws.Range("H95").Value = 48000
$ws.Range("J95").Value = 48000
$ws.Range("L95").Value = 48000
$ws.Range("N95").Value = -53492
# Row 99 (anchor G99=19636)
$ws.Range("H99").Value = 40980.332
$ws.Range("I99").Value = 34499.5
$ws.Range("K99").Value = 34499.5
$ws.Range("M99").Value = -31504.5
# Row 100 (anchor G100=19995)
$ws.Range("H100").Value = 4891
$ws.Range("I100").Value = 3759.4
$ws.Range("K100").Value = 3759.4
$ws.Range("M100").Value = -3218.4
# Row 103 (anchor G103=18526)
$ws.Range("H103").Value = 32197.75
$ws.Range("J103").Value = 32197.75
$ws.Range("L103").Value = 32197.75
$ws.Range("N103").Value = -34541.75
# Row 105 (anchor G105=18698)
$ws.Range("H105").Value = 66897
$ws.Range("J105").Value = 66897
$ws.Range("L105").Value = 66897
$ws.Range("N105").Value = -73885
# Row 106 (anchor G106=18713)
$ws.Range("H106").Value = 33651.2
$ws.Range("J106").Value = 33651.2
$ws.Range("L106").Value = 33651.2
$ws.Range("N106").Value = -36175.2
# Row 110 (anchor G110=25809)
$ws.Range("H110").Value = 42644
$ws.Range("J110").Value = 42644
$ws.Range("L110").Value = 42644
$ws.Range("N110").Value = -50824
# Row 111 (anchor G111=25820)
$ws.Range("H111").Value = 42387
$ws.Range("J111").Value = 42387
$ws.Range("L111").Value = 42387
$ws.Range("N111").Value = -50567
# Row 113 (anchor G113=27740)
$ws.Range("H113").Value = 6842.778
$ws.Range("I113").Value = 4663.8887
$ws.Range("K113").Value = 4663.8887
$ws.Range("M113").Value = -2493.8887
# Row 114 (anchor G114=25990)
$ws.Range("H114").Value = 56958
$ws.Range("J114").Value = 56958
$ws.Range("L114").Value = 56958
$ws.Range("N114").Value = -65636
# Row 132 (anchor G132=44058)
$ws.Range("H132").Value = 9621571
$ws.Range("I132").Value = 17859882
$ws.Range("K132").Value = 53579646
$ws.Range("M132").Value = -53577116
# Row 136 (anchor G136=44060)
$ws.Range("H136").Value = 10825.32
$ws.Range("I136").Value = 3090.4285
$ws.Range("K136").Value = 9271.2855
$ws.Range("M136").Value = -6721.2855

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 46 (anchor G46=42037)
$ws.Range("H46").Value = 90669.28999999999
$ws.Range("J46").Value = 90669.28999999999
$ws.Range("L46").Value = 90669.28999999999
$ws.Range("N46").Value = -91131.28999999999
# Row 122 (anchor G122=36208)
$ws.Range("H122").Value = 186286.69
$ws.Range("I122").Value = 268753.8
$ws.Range("J122").Value = 9571.429
$ws.Range("K122").Value = 806261.3999999999
$ws.Range("L122").Value = 28714.287
$ws.Range("M122").Value = -803811.3999999999
$ws.Range("N122").Value = -33614.287
# Row 125 (anchor G125=34276)
$ws.Range("H125").Value = 71920
$ws.Range("J125").Value = 71920
$ws.Range("L125").Value = 71920
$ws.Range("N125").Value = -81760
# Row 132 (anchor G132=44029)
$ws.Range("H132").Value = 33342396
$ws.Range("I132").Value = 55567204
$ws.Range("K132").Value = 166701612
$ws.Range("M132").Value = -166699082
# Row 134 (anchor G134=42037)
$ws.Range("H134").Value = 90669.28999999999
$ws.Range("J134").Value = 90669.28999999999
$ws.Range("L134").Value = 272007.87
$ws.Range("N134").Value = -277077.87

